$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.136.03"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "2.020.96"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "226.89"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").Value = "0.604"
$ws.Range("E6").Value = "  -2.65%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "55.17"
$ws.Range("E8").Value = "  -3.16%  "
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").Value = "0.0789"
$ws.Range("E10").Value = "  +2.64%  "
$ws.Range("E11").Value = "  -3.13%  "
$ws.Range("D12").Value = "2.321.15"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "14.33"
$ws.Range("E13").Value = "  -2.70%  "
$ws.Range("D14").Value = "20.54"
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("D16").Value = "5.16"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").Value = "2.039.78"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "37.043.11"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "6.21"
$ws.Range("E19").Value = "  +4.20%  "
$ws.Range("D20").Value = "69.04"
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("D21").Value = "0.0₃0826"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").Value = "227.21"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D24").Value = "2.42"
$ws.Range("E25").Value = "  -5.18%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "165.90"
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "9.22"
$ws.Range("E27").Value = "  -2.79%  "
$ws.Range("D28").Value = "0.127"
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("D32").Value = "4.55"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").Value = "0.0619"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("E34").Value = "  -2.71%  "
$ws.Range("E35").Value = "  -4.16%  "
$ws.Range("D36").Value = "1.85"
$ws.Range("E36").Value = "  +1.72%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "3.17"
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("D39").Value = "5.41"
$ws.Range("E39").Value = "  +2.78%  "
$ws.Range("E40").Value = "  -3.83%  "
$ws.Range("D41").Value = "1.481.05"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "16.75"
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("D43").Value = "95.49"
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("D44").Value = "0.0928"
$ws.Range("E44").Value = "  -2.37%  "
$ws.Range("D45").Value = "2.78"
$ws.Range("E45").Value = "  -4.17%  "
$ws.Range("E46").Value = "  -4.10%  "
$ws.Range("D47").Value = "7.33"
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "2.210.55"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").Value = "3.56"
$ws.Range("E51").Value = "  -11.17%  "
